# The document ends with two "Paragraphe de liste" paragraphs:
#   - the second-to-last one carries the hidden _GoBack bookmark
#   - the last one is a bare empty paragraph
#
# The edit expands that tail into six paragraphs:
#   1) empty                               (was the bookmark paragraph, bookmark removed)
#   2) empty                               (new)
#   3) empty                               (new)
#   4) "État de l'art de la Stéganographie :"  (new, ind left=0)
#   5) empty                               (new, ind left=0)
#   6) empty, carries the _GoBack bookmark (was the final paragraph, ind left=0 added)

$d = $word.ActiveDocument

$lastIndex = $d.Paragraphs.Count
$pLast = $d.Paragraphs.Item($lastIndex)
$pBookmark = $d.Paragraphs.Item($lastIndex - 1)
$anchorIndex = $pBookmark.Index

# Move the _GoBack bookmark off of its current paragraph; it will be
# re-added to the (new) final paragraph once the new paragraphs exist.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# Insert four new empty paragraphs, same style, right after $pBookmark
# (i.e. right before $pLast).
$pBookmark.Range.InsertParagraphAfter()
$pBookmark.Range.InsertParagraphAfter()
$pBookmark.Range.InsertParagraphAfter()
$pBookmark.Range.InsertParagraphAfter()

$pHeading = $d.Paragraphs.Item($anchorIndex + 3)
$pBeforeLast = $d.Paragraphs.Item($anchorIndex + 4)
$pFinal = $d.Paragraphs.Item($d.Paragraphs.Count)

# Add the new heading text.
$pHeading.Range.InsertAfter("État de l’art de la Stéganographie :")

# The heading paragraph, the blank paragraph right after it, and the
# final paragraph all get a zero left indent.
$pHeading.Range.ParagraphFormat.LeftIndent = 0
$pBeforeLast.Range.ParagraphFormat.LeftIndent = 0
$pFinal.Range.ParagraphFormat.LeftIndent = 0

# Re-create the _GoBack bookmark on the (new) final paragraph.
$d.Bookmarks.Add("_GoBack", $pFinal.Range)
